# Applies the edit described by the diff:
# - Populates column I ("Email") with per-row email addresses for rows 2-106
#   of the first worksheet (Hoja1), skipping the blank placeholder rows
#   (9, 28, 59, 79, 100) that have no other data.
# - Populates column L ("IBAN") for the 4 rows that were previously blank
#   (rows 2, 3, 13, 17) with their IBAN values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$emails = @{
    2 = 'am00@tecnoleonsl.es'
    3 = 'aam00@tecnoproyectsl.es'
    4 = 'aab00@asoftware.es'
    5 = 'acd00@asoftware.es'
    6 = 'abg00@tecnoleonsl.es'
    7 = 'agm00@tecnoproyectsl.es'
    8 = 'ari00@tecnoleonsl.es'
    10 = 'amg00@asoftware.es'
    11 = 'afr00@pblanksa.es'
    12 = 'ala00@tecnoproyectsl.es'
    13 = 'ala01@tecnoproyectsl.es'
    14 = 'agr00@aphonsa.es'
    15 = 'ads00@aphonsa.es'
    16 = 'bv00@tecnoproyectsl.es'
    17 = 'bbv00@tecnoproyectsl.es'
    18 = 'bba00@pblanksa.es'
    19 = 'bam00@pblanksa.es'
    20 = 'bca00@pblanksa.es'
    21 = 'bbl00@pblanksa.es'
    22 = 'bga00@tecnoproyectsl.es'
    23 = 'bl00@asoftware.es'
    24 = 'bbl00@tecnoleonsl.es'
    25 = 'bs00@tecnoproyectsl.es'
    26 = 'bga00@tecnoleonsl.es'
    27 = 'bgd00@asoftware.es'
    29 = 'cim00@tecnoproyectsl.es'
    30 = 'cqc00@asoftware.es'
    31 = 'cvm00@pblanksa.es'
    32 = 'cpf00@tecnoproyectsl.es'
    33 = 'cci00@pblanksa.es'
    34 = 'cfs00@tecnoleonsl.es'
    35 = 'cje00@tecnoleonsl.es'
    36 = 'cgg00@aphonsa.es'
    37 = 'cpe00@asoftware.es'
    38 = 'csm00@pblanksa.es'
    39 = 'clm00@tecnoleonsl.es'
    40 = 'ccj00@pblanksa.es'
    41 = 'cpj00@tecnoproyectsl.es'
    42 = 'ccb00@tecnoproyectsl.es'
    43 = 'can00@pblanksa.es'
    44 = 'ccm00@tecnoproyectsl.es'
    45 = 'ccm00@tecnoproyectsl.es'
    46 = 'dam00@tecnoproyectsl.es'
    47 = 'dfc00@aphonsa.es'
    48 = 'dlc00@asoftware.es'
    49 = 'eah00@tecnoleonsl.es'
    50 = 'ek00@tecnoproyectsl.es'
    51 = 'es00@asoftware.es'
    52 = 'fpb00@tecnoproyectsl.es'
    53 = 'faa00@tecnoleonsl.es'
    54 = 'ffi00@tecnoleonsl.es'
    55 = 'ffd00@aphonsa.es'
    56 = 'fme00@aphonsa.es'
    57 = 'flc00@aphonsa.es'
    58 = 'gar00@tecnoleonsl.es'
    60 = 'gpl00@tecnoleonsl.es'
    61 = 'gge00@asoftware.es'
    62 = 'gpl00@tecnoproyectsl.es'
    63 = 'gfd00@aphonsa.es'
    64 = 'gsa00@asoftware.es'
    65 = 'gif00@tecnoleonsl.es'
    66 = 'gac00@aphonsa.es'
    67 = 'gra00@tecnoproyectsl.es'
    68 = 'gpl00@tecnoproyectsl.es'
    69 = 'gmn00@pblanksa.es'
    70 = 'glj00@aphonsa.es'
    71 = 'hdd00@asoftware.es'
    72 = 'ldr00@pblanksa.es'
    73 = 'lrf00@tecnoproyectsl.es'
    74 = 'mcn00@tecnoproyectsl.es'
    75 = 'mbo00@aphonsa.es'
    76 = 'mfj00@tecnoproyectsl.es'
    77 = 'mgj00@aphonsa.es'
    78 = 'moh00@tecnoleonsl.es'
    80 = 'mvp00@tecnoproyectsl.es'
    81 = 'mmc00@aphonsa.es'
    82 = 'mbm00@tecnoproyectsl.es'
    83 = 'nbb00@asoftware.es'
    84 = 'osp00@aphonsa.es'
    85 = 'odp00@pblanksa.es'
    86 = 'ofj00@aphonsa.es'
    87 = 'pdc00@asoftware.es'
    88 = 'pgc00@asoftware.es'
    89 = 'ppv00@tecnoleonsl.es'
    90 = 'pac00@asoftware.es'
    91 = 'pma00@asoftware.es'
    92 = 'pgr00@pblanksa.es'
    93 = 'rfi00@tecnoleonsl.es'
    94 = 'raf00@tecnoleonsl.es'
    95 = 'rdf00@tecnoleonsl.es'
    96 = 'rgt00@tecnoleonsl.es'
    97 = 'rlc00@aphonsa.es'
    98 = 'sbc00@asoftware.es'
    99 = 'sas00@asoftware.es'
    101 = 'she00@asoftware.es'
    102 = 'sld00@asoftware.es'
    103 = 'sfh00@tecnoleonsl.es'
    104 = 'sgs00@aphonsa.es'
    105 = 'ssg00@tecnoleonsl.es'
    106 = 'sss00@tecnoleonsl.es'
}

$ibans = @{
    2 = 'ES3520960043043554600000'
    3 = 'ES7832145464138452163421'
    13 = 'ES7921564975243245467995'
    17 = 'ES0721584976902154655487'
}

foreach ($row in $emails.Keys) {
    $ws.Cells.Item([int]$row, 9).Value = $emails[$row]
}

foreach ($row in $ibans.Keys) {
    $ws.Cells.Item([int]$row, 12).Value = $ibans[$row]
}
